$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 589
$ws.Range("I2").Value = 297
$ws.Range("J2").Value = 735
$ws.Range("K2").Value = 297
$ws.Range("L2").Value = 735
$ws.Range("M2").Value = -184
$ws.Range("N2").Value = -961

$ws.Range("H38").Value = 4321.8887
$ws.Range("I38").Value = 2699.5715
$ws.Range("J38").Value = 10000
$ws.Range("K38").Value = 8098.7145
$ws.Range("L38").Value = 30000
$ws.Range("M38").Value = -7726.7145
$ws.Range("N38").Value = -30744

$ws.Range("H40").Value = 17254668
$ws.Range("J40").Value = 35724690
$ws.Range("L40").Value = 35724690
$ws.Range("N40").Value = -35725040

$ws.Range("H43").Value = 5630204.5
$ws.Range("I43").Value = 8442308
$ws.Range("J43").Value = 5999
$ws.Range("K43").Value = 8442308
$ws.Range("L43").Value = 5999
$ws.Range("M43").Value = -8442239
$ws.Range("N43").Value = -6137

$ws.Range("H58").Value = 2219.8
$ws.Range("I58").Value = 350
$ws.Range("J58").Value = 5024.5
$ws.Range("K58").Value = 1050
$ws.Range("L58").Value = 15073.5
$ws.Range("M58").Value = -900
$ws.Range("N58").Value = -15373.5

$ws.Range("H113").Value = 3425.75
$ws.Range("I113").Value = 2822.3
$ws.Range("J113").Value = 3856.7856
$ws.Range("K113").Value = 2822.3
$ws.Range("L113").Value = 3856.7856
$ws.Range("M113").Value = 431.6999999999998
$ws.Range("N113").Value = -10364.7856

$ws.Range("H116").Value = 41673136
$ws.Range("I116").Value = 24627656
$ws.Range("J116").Value = 166673330
$ws.Range("K116").Value = 24627656
$ws.Range("L116").Value = 166673330
$ws.Range("M116").Value = -24624214
$ws.Range("N116").Value = -166680214

$ws.Range("H137").Value = 19078046
$ws.Range("I137").Value = 1111821.1
$ws.Range("J137").Value = 37044268
$ws.Range("K137").Value = 3335463.3
$ws.Range("L137").Value = 111132804
$ws.Range("M137").Value = -3332913.3
$ws.Range("N137").Value = -111137904

$ws.Range("H138").Value = 3884.2222
$ws.Range("I138").Value = 2838.7273
$ws.Range("J138").Value = 5527.143
$ws.Range("K138").Value = 8516.1819
$ws.Range("L138").Value = 16581.429
$ws.Range("M138").Value = -3376.1819
$ws.Range("N138").Value = -26861.429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 5817.3335
$ws.Range("I122").Value = 4734.3335
$ws.Range("K122").Value = 14203.0005
$ws.Range("M122").Value = -11753.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5382.84
$ws.Range("I20").Value = 4877.8237
$ws.Range("K20").Value = 4877.8237
$ws.Range("M20").Value = -4630.8237

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 7891.8335
$ws.Range("I122").Value = 1675.5
$ws.Range("K122").Value = 5026.5
$ws.Range("M122").Value = -2576.5

$ws.Range("H134").Value = 2172.2
$ws.Range("J134").Value = 3730.6667
$ws.Range("L134").Value = 11192.0001
$ws.Range("N134").Value = -16262.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 2287.2083
$ws.Range("J55").Value = 3039.2307
$ws.Range("L55").Value = 9117.6921
$ws.Range("N55").Value = -9471.6921

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 2846752.8
$ws.Range("I70").Value = 3501903.2
$ws.Range("K70").Value = 3501903.2
$ws.Range("M70").Value = -3501633.2

$ws.Range("H73").Value = 2846752.8
$ws.Range("I73").Value = 3501903.2
$ws.Range("K73").Value = 3501903.2
$ws.Range("M73").Value = -3500967.2

$ws.Range("H102").Value = 17248450
$ws.Range("I102").Value = 23816708
$ws.Range("K102").Value = 23816708
$ws.Range("M102").Value = -23815086

$ws.Range("H113").Value = 2009.7778
$ws.Range("I113").Value = 1797.8
$ws.Range("J113").Value = 2274.75
$ws.Range("K113").Value = 1797.8
$ws.Range("L113").Value = 2274.75
$ws.Range("M113").Value = 372.2
$ws.Range("N113").Value = -6614.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 973.1429000000001
$ws.Range("I22").Value = 990.5
$ws.Range("J22").Value = 950
$ws.Range("K22").Value = 990.5
$ws.Range("L22").Value = 950
$ws.Range("M22").Value = -695.5
$ws.Range("N22").Value = -1540

$ws.Range("H27").Value = 973.1429000000001
$ws.Range("I27").Value = 990.5
$ws.Range("J27").Value = 950
$ws.Range("K27").Value = 990.5
$ws.Range("L27").Value = 950
$ws.Range("M27").Value = -883.5
$ws.Range("N27").Value = -1164

$ws.Range("H40").Value = 18523836
$ws.Range("I40").Value = 2599.875
$ws.Range("J40").Value = 33340824
$ws.Range("K40").Value = 2599.875
$ws.Range("L40").Value = 33340824
$ws.Range("M40").Value = -2463.875
$ws.Range("N40").Value = -33341096

$ws.Range("H46").Value = 6937.0312
$ws.Range("I46").Value = 5401.25
$ws.Range("J46").Value = 7156.4287
$ws.Range("K46").Value = 5401.25
$ws.Range("L46").Value = 7156.4287
$ws.Range("M46").Value = -5213.25
$ws.Range("N46").Value = -7532.4287

$ws.Range("H55").Value = 374.6
$ws.Range("I55").Value = 135.66667
$ws.Range("J55").Value = 533.8889
$ws.Range("K55").Value = 135.66667
$ws.Range("L55").Value = 533.8889
$ws.Range("M55").Value = 37.33332999999999
$ws.Range("N55").Value = -879.8889

$ws.Range("H100").Value = 5908.722
$ws.Range("I100").Value = 4057.9092
$ws.Range("J100").Value = 8817.143
$ws.Range("K100").Value = 4057.9092
$ws.Range("L100").Value = 8817.143
$ws.Range("M100").Value = -3516.9092
$ws.Range("N100").Value = -9899.143

$ws.Range("H122").Value = 40005770
$ws.Range("I122").Value = 66671284
$ws.Range("K122").Value = 200013852
$ws.Range("M122").Value = -200011402

$ws.Range("H132").Value = 4925.2915
$ws.Range("I132").Value = 4583.3125
$ws.Range("J132").Value = 5609.25
$ws.Range("K132").Value = 13749.9375
$ws.Range("L132").Value = 16827.75
$ws.Range("M132").Value = -11219.9375
$ws.Range("N132").Value = -21887.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5284.143
$ws.Range("J62").Value = 8399.799999999999
$ws.Range("L62").Value = 8399.799999999999
$ws.Range("N62").Value = -9647.799999999999

$ws.Range("H65").Value = 5284.143
$ws.Range("J65").Value = 8399.799999999999
$ws.Range("L65").Value = 41999
$ws.Range("N65").Value = -48239

$ws.Range("H81").Value = 6944824
$ws.Range("I81").Value = 10417038
$ws.Range("K81").Value = 20834076
$ws.Range("M81").Value = -20833015

$ws.Range("H84").Value = 6944824
$ws.Range("I84").Value = 10417038
$ws.Range("K84").Value = 104170380
$ws.Range("M84").Value = -104165076

$ws.Range("H107").Value = 2493.0312
$ws.Range("I107").Value = 2558.0454
$ws.Range("K107").Value = 7674.1362
$ws.Range("M107").Value = -5754.1362
